# Applies the "weight column" model update to GDP_SDMX_Model.xlsx
#  * GDP_TYPE sheet/concept is repurposed into a TRANSFORMATION sheet/concept
#  * INDICATOR sheet gains new weight / net-tax indicator rows
#  * the old GDP_TYPE codelist (Norminal/Real GDP) is replaced by the new
#    TRANSFORMATION codelist (PM1, YM1, EST, ICLB, ICUB)
#  * INDUSTRY_TYPE drops the now-obsolete "Net Taxes" / "GDP" rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. DSD sheet: the concept row that used to describe GDP_TYPE now
#    describes TRANSFORMATION
# ---------------------------------------------------------------------
$dsd = $wb.Worksheets.Item("DSD")
$dsd.Range("A6").Value = "TRANSFORMATION"
$dsd.Range("B6").Value = "Transformation"
$dsd.Range("F6").Value = "CL_GDP_TYPE"
$dsd.Range("C6").Select()

# ---------------------------------------------------------------------
# 2. INDICATOR sheet: rework the GDP rows + add the new weight/net-tax
#    indicator rows (this is the "weight column" generation the commit
#    message refers to)
# ---------------------------------------------------------------------
$indicator = $wb.Worksheets.Item("INDICATOR")
$indicator.Range("A2").Value = "NRGDP"
$indicator.Range("B2").Value = "Norminal GDP"
$indicator.Range("A3").Value = "RLGDP"
$indicator.Range("B3").Value = "Real GDP"
$indicator.Range("A4").Value = "NRWGT"
$indicator.Range("B4").Value = "Norminal GDP weight"
$indicator.Range("A5").Value = "RLWGT"
$indicator.Range("B5").Value = "Real GDP weight"
$indicator.Range("A6").Value = "NRTAX"
$indicator.Range("B6").Value = "Norminal GDP net tax"
$indicator.Range("A7").Value = "RLTAX"
$indicator.Range("B7").Value = "Real GDP net tax"
$indicator.Columns.Item(1).AutoFit()
$indicator.Columns.Item(2).AutoFit()

# ---------------------------------------------------------------------
# 3. GDP_TYPE sheet becomes the TRANSFORMATION codelist sheet: rename
#    the tab, replace the Norminal/Real GDP rows with the new
#    transformation codes and append the estimate / IC-bound codes
# ---------------------------------------------------------------------
$gdpType = $wb.Worksheets.Item("GDP_TYPE")
$gdpType.Range("A2").Value = "PM1"
$gdpType.Range("B2").Value = "% change over previous period"
$gdpType.Range("A3").Value = "YM1"
$gdpType.Range("B3").Value = "% change over same period last year"
$gdpType.Range("B2:B3").Font.Name = "Arial"
$gdpType.Range("B2:B3").Font.Size = 9
$gdpType.Range("A4").Value = "EST"
$gdpType.Range("B4").Value = "Estimated value"
$gdpType.Range("A5").Value = "ICLB"
$gdpType.Range("B5").Value = "IC lower bound"
$gdpType.Range("A6").Value = "ICUB"
$gdpType.Range("B6").Value = "IC upper bound"
$gdpType.Name = "TRANSFORMATION"

# ---------------------------------------------------------------------
# 4. INDUSTRY_TYPE sheet: the "Net Taxes" / "GDP" rows are no longer
#    needed now that GDP_TYPE has become TRANSFORMATION
# ---------------------------------------------------------------------
$industry = $wb.Worksheets.Item("INDUSTRY_TYPE")
$industry.Rows("24:25").Delete()
$industry.Range("B25").Select()
